$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lugar")

# Append rows 3-5 first (new blind/point data)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 27.35455
$ws.Range("C3").Value = -99.937733333333298
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Blind Los Conejos"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 27.3187444
$ws.Range("C4").Value = -99.966144444444396
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Blind Diego"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 27.3271278
$ws.Range("C5").Value = -99.968694444444395
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "Rancho El Huisachito"

# Update row 2 - the "No recuerdo el Nombre" entry gets corrected/updated
$ws.Range("B2").Value = 27.319050000000001
$ws.Range("C2").Value = -99.989744166666597
$ws.Range("E2").Value = "Blind La Cuchillita"

# Append remaining rows 6-8
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 27.332599999999999
$ws.Range("C6").Value = -99.958705555555497
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Blind Puertas Amarillas"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 27.335111099999999
$ws.Range("C7").Value = -99.966255555555506
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Blind Cuatro caminos"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 27.3191278
$ws.Range("C8").Value = -99.986011111111097
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Blind La Cuchilla"

# Move selection to A9, matching the post-edit cursor position
$ws.Range("A9").Select()
